$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the newly-completed data for row 11 (Issue #9 is now fixed)
$ws.Range("C11").Value2 = "FIX"
$ws.Range("E11").Value2 = "1.3.3"
$ws.Range("J11").Value2 = "Functionality restored."

# Match row 11's formatting (fill/alignment) to the already-"fixed" row above it (row 10)
$ws.Range("A10:J10").Copy()
$ws.Range("A11:J11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Scroll the sheet view so row 7 is at the top-left, as in the saved workbook
[void]$ws.Range("A7").Select()
$excel.ActiveWindow.ScrollRow = 7
